$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.451.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.16%  '

# Row 3
$ws.Range("D3").Value = "'3.759.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.39%  '

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = "'616.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.77%  '

# Row 6
$ws.Range("D6").Value = "'179.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.04%  '

# Row 7
$ws.Range("D7").Value = "'3.755.04"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.34%  '

# Row 8
$ws.Range("E8").Value = '  +0.09%  '

# Row 9
$ws.Range("E9").Value = '  -1.10%  '

# Row 10
$ws.Range("E10").Value = '  +1.17%  '

# Row 11
$ws.Range("D11").Value = "'6.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.30%  '

# Row 12
$ws.Range("E12").Value = '  -1.79%  '

# Row 13
$ws.Range("D13").Value = "'40.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.20%  '

# Row 14
$ws.Range("E14").Value = '  +1.26%  '

# Row 15
$ws.Range("D15").Value = "'4.385.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.41%  '

# Row 16
$ws.Range("D16").Value = "'3.758.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.32%  '

# Row 17
$ws.Range("D17").Value = "'69.564.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.12%  '

# Row 18
$ws.Range("E18").Value = '  -2.44%  '

# Row 19
$ws.Range("E19").Value = '  -1.35%  '

# Row 20
$ws.Range("D20").Value = "'16.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.25%  '

# Row 21
$ws.Range("D21").Value = "'501.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.14%  '

# Row 22
$ws.Range("D22").Value = "'9.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.26%  '

# Row 23
$ws.Range("D23").Value = "'0.722"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '

# Row 24
$ws.Range("E24").Value = '  +3.18%  '

# Row 25
$ws.Range("D25").Value = "'85.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.18%  '

# Row 26
$ws.Range("E26").Value = '  -2.00%  '

# Row 27
$ws.Range("D27").Value = "'10.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.64%  '

# Row 28
$ws.Range("E28").Value = '  +6.58%  '

# Row 29
$ws.Range("E29").Value = '  +0.09%  '

# Row 30
$ws.Range("D30").Value = "'2.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.59%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'2.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.13%  '

# Row 32
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = "'8.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.15%  '

# Row 33
$ws.Range("D33").Value = "'30.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.63%  '

# Row 34
$ws.Range("E34").Value = '  -1.00%  '

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.07%  '

# Row 36
$ws.Range("E36").Value = '  +0.94%  '

# Row 37
$ws.Range("D37").Value = "'6.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.71%  '

# Row 38
$ws.Range("E38").Value = '  +4.28%  '

# Row 39
$ws.Range("E39").Value = '  +4.21%  '

# Row 40
$ws.Range("D40").Value = "'466.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.75%  '

# Row 41
$ws.Range("D41").Value = "'3.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +13.73%  '

# Row 42
$ws.Range("E42").Value = '  -4.13%  '

# Row 43
$ws.Range("D43").Value = "'49.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.77%  '

# Row 44
$ws.Range("D44").Value = "'44.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.54%  '

# Row 45
$ws.Range("D45").Value = "'8.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.00%  '

# Row 46
$ws.Range("D46").Value = "'2.957.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.31%  '

# Row 47
$ws.Range("E47").Value = '  -0.20%  '

# Row 48
$ws.Range("D48").Value = "'27.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.39%  '

# Row 49
$ws.Range("D49").Value = "'139.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.22%  '

# Row 50
$ws.Range("E50").Value = '  +0.01%  '

# Row 51
$ws.Range("D51").Value = "'2.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.31%  '
